$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A88").Value = 87
$ws.Range("B88").Value = 1
$ws.Range("C88").Value = "2024-06-16 18:16:24"
$ws.Range("D88").Value = 200
$ws.Range("E88").Value = 8

$ws.Range("A89").Value = 88
$ws.Range("B89").Value = 2
$ws.Range("C89").Value = "2024-06-16 18:16:24"
$ws.Range("D89").Value = 200
$ws.Range("E89").Value = 0
